$p = $ppt.ActivePresentation

# --- Slide 1 : title-slide SmartArt ("Guided By - Mr. Manoj Pawaiya Sir") ---
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(4)
$sa1 = $sh1.SmartArt
$n1 = $sa1.AllNodes.Item(1)
$n1.TextFrame2.TextRange.Text = "Guided By " + [char]0x2013 + " Manoj Pawaiya Sir"

# --- Slide 12 : skills SmartArt ("HTML, CSS, JavaScript,JSON (Front-End)") ---
$s12 = $p.Slides.Item(12)
$sh12 = $s12.Shapes.Item(2)
$sa12 = $sh12.SmartArt
$n12 = $sa12.AllNodes.Item(3)
$n12.TextFrame2.TextRange.Text = "HTML, CSS, JavaScript, JSON (Front-End)"
